$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '74.902.39'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.823.47'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +7.54%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.68'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '595.46'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.551'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.09%  '
$ws.Range("E9").Value = '  -4.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.821.79'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +7.64%  '
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.371'
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.89'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.340.59'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.729.06'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000187'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.80'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.820.76'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +7.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.94'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.32'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.47'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.91'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.965.24'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +7.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.16'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.74'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.50%  '
$ws.Range("E29").Value = '  +10.13%  '
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '518.34'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.73'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.79'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.64%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.27'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.94'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.98%  '
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '186.37'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +15.54%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.341'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.99'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.59%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  +1.90%  '
$ws.Range("E46").Value = '  +2.70%  '
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0854'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.579'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +9.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.71'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("E51").Value = '  +8.28%  '

Write-Output "Applied 80 cell updates"